$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.101729968414588
$ws.Cells.Item(2, 3).Value = 0.3226985537384337
$ws.Cells.Item(2, 4).Value = 0.2735004925934419
$ws.Cells.Item(2, 6).Value = 1.484329595726805
$ws.Cells.Item(2, 7).Value = 0.8165580542823108
$ws.Cells.Item(2, 8).Value = 0.8445637225456437
$ws.Cells.Item(2, 9).Value = 0.5751832815021061
$ws.Cells.Item(2, 10).Value = 0.3485596430322886
$ws.Cells.Item(2, 14).Value = 1.02719948783637

$ws.Cells.Item(3, 2).Value = 0.9902222696572949
$ws.Cells.Item(3, 3).Value = 0.2859861336109759
$ws.Cells.Item(3, 4).Value = 0.2674005364426364
$ws.Cells.Item(3, 6).Value = 1.467441684983953
$ws.Cells.Item(3, 7).Value = 0.8028150664091243
$ws.Cells.Item(3, 8).Value = 0.8439880063391456
$ws.Cells.Item(3, 9).Value = 0.5781526834932151
$ws.Cells.Item(3, 10).Value = 0.3371667947527044
$ws.Cells.Item(3, 14).Value = 1.039502143921151

$ws.Cells.Item(4, 2).Value = 0.9219430868801624
$ws.Cells.Item(4, 3).Value = 0.2634805851665192
$ws.Cells.Item(4, 4).Value = 0.2637629761289304
$ws.Cells.Item(4, 6).Value = 1.458135037726933
$ws.Cells.Item(4, 7).Value = 0.7950871167749227
$ws.Cells.Item(4, 8).Value = 0.8441741054577676
$ws.Cells.Item(4, 9).Value = 0.5804232672984675
$ws.Cells.Item(4, 10).Value = 0.3303880654075186
$ws.Cells.Item(4, 14).Value = 1.047532604180645

$ws.Cells.Item(5, 2).Value = 0.8941664184201841
$ws.Cells.Item(5, 3).Value = 0.2543184456173719
$ws.Cells.Item(5, 4).Value = 0.2623078126727165
$ws.Cells.Item(5, 6).Value = 1.454608854962544
$ws.Cells.Item(5, 7).Value = 0.7921156237053282
$ws.Cells.Item(5, 8).Value = 0.8443853366069476
$ws.Cells.Item(5, 9).Value = 0.5814606556403561
$ws.Cells.Item(5, 10).Value = 0.3276799150211502
$ws.Cells.Item(5, 14).Value = 1.050924881274813

$ws.Cells.Item(6, 2).Value = 0.8895570208416075
$ws.Cells.Item(6, 3).Value = 0.2527976261742708
$ws.Cells.Item(6, 4).Value = 0.2620678266611094
$ws.Cells.Item(6, 6).Value = 1.454039397328174
$ws.Cells.Item(6, 7).Value = 0.7916329164419977
$ws.Cells.Item(6, 8).Value = 0.8444285804205407
$ws.Cells.Item(6, 9).Value = 0.5816396729329583
$ws.Cells.Item(6, 10).Value = 0.3272334983921752
$ws.Cells.Item(6, 14).Value = 1.051495400821256

$ws.Cells.Item(7, 2).Value = 0.9215682872265347
$ws.Cells.Item(7, 3).Value = 0.2633569846887553
$ws.Cells.Item(7, 4).Value = 0.2637432411925005
$ws.Cells.Item(7, 6).Value = 1.458086405072549
$ws.Cells.Item(7, 7).Value = 0.7950463238387897
$ws.Cells.Item(7, 8).Value = 0.8441764063784518
$ws.Cells.Item(7, 9).Value = 0.5804368044772019
$ws.Cells.Item(7, 10).Value = 0.3303513230157051
$ws.Cells.Item(7, 14).Value = 1.047577868608492

$ws.Cells.Item(8, 2).Value = 1.063243643365865
$ws.Cells.Item(8, 3).Value = 0.3100326186583686
$ws.Cells.Item(8, 4).Value = 0.2713748765025628
$ws.Cells.Item(8, 6).Value = 1.478285580746416
$ws.Cells.Item(8, 7).Value = 0.8116714989894405
$ws.Cells.Item(8, 8).Value = 0.8442530453320813
$ws.Cells.Item(8, 9).Value = 0.5761140488638361
$ws.Cells.Item(8, 10).Value = 0.3445863176967805
$ws.Cells.Item(8, 14).Value = 1.031342504390274

$ws.Cells.Item(9, 2).Value = 1.342539431794535
$ws.Cells.Item(9, 3).Value = 0.4018543644489228
$ws.Cells.Item(9, 4).Value = 0.2871947701663657
$ws.Cells.Item(9, 6).Value = 1.526369994461476
$ws.Cells.Item(9, 7).Value = 0.8499530057626714
$ws.Cells.Item(9, 8).Value = 0.8486991842206351
$ws.Cells.Item(9, 9).Value = 0.571204639797493
$ws.Cells.Item(9, 10).Value = 0.3742305264130863
$ws.Cells.Item(9, 14).Value = 1.003288247185537

$ws.Cells.Item(10, 2).Value = 1.548639413767887
$ws.Cells.Item(10, 3).Value = 0.4695073440162219
$ws.Cells.Item(10, 4).Value = 0.2993381918975473
$ws.Cells.Item(10, 6).Value = 1.566928962198574
$ws.Cells.Item(10, 7).Value = 0.8816056427920671
$ws.Cells.Item(10, 8).Value = 0.8546065499545819
$ws.Cells.Item(10, 9).Value = 0.5697981136392301
$ws.Cells.Item(10, 10).Value = 0.397082964746815
$ws.Cells.Item(10, 14).Value = 0.9849855957703753

$ws.Cells.Item(11, 2).Value = 1.642597921449067
$ws.Cells.Item(11, 3).Value = 0.5003295612153806
$ws.Cells.Item(11, 4).Value = 0.3049756481985355
$ws.Cells.Item(11, 6).Value = 1.586530820022247
$ws.Cells.Item(11, 7).Value = 0.8967854275081208
$ws.Cells.Item(11, 8).Value = 0.8578722846403082
$ws.Cells.Item(11, 9).Value = 0.5696416417728685
$ws.Cells.Item(11, 10).Value = 0.4077162794485076
$ws.Cells.Item(11, 14).Value = 0.9771612401210703

$ws.Cells.Item(12, 2).Value = 1.678206434838103
$ws.Cells.Item(12, 3).Value = 0.5120079625346534
$ws.Cells.Item(12, 4).Value = 0.307126680047304
$ws.Cells.Item(12, 6).Value = 1.594120171200302
$ws.Cells.Item(12, 7).Value = 0.9026469793456613
$ws.Cells.Item(12, 8).Value = 0.8591924831536915
$ws.Cells.Item(12, 9).Value = 0.5696523444796568
$ws.Cells.Item(12, 10).Value = 0.4117772952157566
$ws.Cells.Item(12, 14).Value = 0.9742705761679957

$ws.Cells.Item(13, 2).Value = 1.670536248818223
$ws.Cells.Item(13, 3).Value = 0.5094925119031473
$ws.Cells.Item(13, 4).Value = 0.3066626949430429
$ws.Cells.Item(13, 6).Value = 1.592478245350961
$ws.Cells.Item(13, 7).Value = 0.9013795341849402
$ws.Cells.Item(13, 8).Value = 0.8589044338005749
$ws.Cells.Item(13, 9).Value = 0.5696469210197321
$ws.Cells.Item(13, 10).Value = 0.4109011485536485
$ws.Cells.Item(13, 14).Value = 0.9748899170824714

$ws.Cells.Item(14, 2).Value = 1.6455268885544
$ws.Cells.Item(14, 3).Value = 0.5012902160689237
$ws.Cells.Item(14, 4).Value = 0.3051522895039227
$ws.Cells.Item(14, 6).Value = 1.587151856863628
$ws.Cells.Item(14, 7).Value = 0.8972653844471381
$ws.Cells.Item(14, 8).Value = 0.8579792219350679
$ws.Cells.Item(14, 9).Value = 0.5696411177089473
$ws.Cells.Item(14, 10).Value = 0.4080496906424003
$ws.Cells.Item(14, 14).Value = 0.976921974762611

$ws.Cells.Item(15, 2).Value = 1.630211620769785
$ws.Cells.Item(15, 3).Value = 0.4962669453414605
$ws.Cells.Item(15, 4).Value = 0.3042292372467443
$ws.Cells.Item(15, 6).Value = 1.583911009643757
$ws.Cells.Item(15, 7).Value = 0.8947601336838602
$ws.Cells.Item(15, 8).Value = 0.8574233923905297
$ws.Cells.Item(15, 9).Value = 0.569646686588257
$ws.Cells.Item(15, 10).Value = 0.4063075787382076
$ws.Cells.Item(15, 14).Value = 0.9781760814476641

$ws.Cells.Item(16, 2).Value = 1.542503023693371
$ws.Cells.Item(16, 3).Value = 0.4674939763743851
$ws.Cells.Item(16, 4).Value = 0.298972047474777
$ws.Cells.Item(16, 6).Value = 1.5656711874975
$ws.Cells.Item(16, 7).Value = 0.8806294071304137
$ws.Cells.Item(16, 8).Value = 0.85440479664598
$ws.Cells.Item(16, 9).Value = 0.5698181042039252
$ws.Cells.Item(16, 10).Value = 0.3963928547140085
$ws.Cells.Item(16, 14).Value = 0.9855070395028918

$ws.Cells.Item(17, 2).Value = 1.488748050397987
$ws.Cells.Item(17, 3).Value = 0.4498546130399177
$ws.Cells.Item(17, 4).Value = 0.2957759278976511
$ws.Cells.Item(17, 6).Value = 1.554777242381249
$ws.Cells.Item(17, 7).Value = 0.8721613750870461
$ws.Cells.Item(17, 8).Value = 0.8527013915969235
$ws.Cells.Item(17, 9).Value = 0.5700473824343604
$ws.Cells.Item(17, 10).Value = 0.3903715242232124
$ws.Cells.Item(17, 14).Value = 0.9901329018281189

$ws.Cells.Item(18, 2).Value = 1.457848719491665
$ws.Cells.Item(18, 3).Value = 0.4397132975601608
$ws.Cells.Item(18, 4).Value = 0.2939482767028352
$ws.Cells.Item(18, 6).Value = 1.54861964719629
$ws.Cells.Item(18, 7).Value = 0.8673642226362119
$ws.Cells.Item(18, 8).Value = 0.8517760624939115
$ws.Cells.Item(18, 9).Value = 0.5702247267779015
$ws.Cells.Item(18, 10).Value = 0.3869305661481093
$ws.Cells.Item(18, 14).Value = 0.9928407820974101

$ws.Cells.Item(19, 2).Value = 1.447390039317042
$ws.Cells.Item(19, 3).Value = 0.4362803745587485
$ws.Cells.Item(19, 4).Value = 0.2933312999399078
$ws.Cells.Item(19, 6).Value = 1.546553364366616
$ws.Cells.Item(19, 7).Value = 0.8657525711068388
$ws.Cells.Item(19, 8).Value = 0.8514720986100315
$ws.Cells.Item(19, 9).Value = 0.5702925683520164
$ws.Cells.Item(19, 10).Value = 0.385769348572822
$ws.Cells.Item(19, 14).Value = 0.9937657281068368

$ws.Cells.Item(20, 2).Value = 1.494468386339747
$ws.Cells.Item(20, 3).Value = 0.4517319000647149
$ws.Cells.Item(20, 4).Value = 0.2961150561143313
$ws.Cells.Item(20, 6).Value = 1.555925704266514
$ws.Cells.Item(20, 7).Value = 0.8730552039435793
$ws.Cells.Item(20, 8).Value = 0.8528770867015112
$ws.Cells.Item(20, 9).Value = 0.5700182662235207
$ws.Cells.Item(20, 10).Value = 0.3910101900573295
$ws.Cells.Item(20, 14).Value = 0.9896355843625813

$ws.Cells.Item(21, 2).Value = 1.652871976597851
$ws.Cells.Item(21, 3).Value = 0.5036992473790747
$ws.Cells.Item(21, 4).Value = 0.3055954915841141
$ws.Cells.Item(21, 6).Value = 1.588711818174843
$ws.Cells.Item(21, 7).Value = 0.8984707270105048
$ws.Cells.Item(21, 8).Value = 0.8582487090969266
$ws.Cells.Item(21, 9).Value = 0.5696409201336436
$ws.Cells.Item(21, 10).Value = 0.4088862971411231
$ws.Cells.Item(21, 14).Value = 0.9763231483481647

$ws.Cells.Item(22, 2).Value = 1.756563499613776
$ws.Cells.Item(22, 3).Value = 0.5377019512721972
$ws.Cells.Item(22, 4).Value = 0.3118861961509367
$ws.Cells.Item(22, 6).Value = 1.611110700984781
$ws.Cells.Item(22, 7).Value = 0.9157420952181496
$ws.Cells.Item(22, 8).Value = 0.8622464525027738
$ws.Cells.Item(22, 9).Value = 0.5698022649393195
$ws.Cells.Item(22, 10).Value = 0.4207700491824085
$ws.Cells.Item(22, 14).Value = 0.968043888440782

$ws.Cells.Item(23, 2).Value = 1.701206444988941
$ws.Cells.Item(23, 3).Value = 0.5195504935526287
$ws.Cells.Item(23, 4).Value = 0.3085200807819319
$ws.Cells.Item(23, 6).Value = 1.599066797204515
$ws.Cells.Item(23, 7).Value = 0.9064632365795262
$ws.Cells.Item(23, 8).Value = 0.8600680953794324
$ws.Cells.Item(23, 9).Value = 0.569678674884166
$ws.Cells.Item(23, 10).Value = 0.4144090245785037
$ws.Cells.Item(23, 14).Value = 0.972424105719746

$ws.Cells.Item(24, 2).Value = 1.491882206464709
$ws.Cells.Item(24, 3).Value = 0.4508831792612114
$ws.Cells.Item(24, 4).Value = 0.2959617055928589
$ws.Cells.Item(24, 6).Value = 1.55540615616701
$ws.Cells.Item(24, 7).Value = 0.87265088211538
$ws.Cells.Item(24, 8).Value = 0.8527974868233628
$ws.Cells.Item(24, 9).Value = 0.5700312878764393
$ws.Cells.Item(24, 10).Value = 0.390721384512787
$ws.Cells.Item(24, 14).Value = 0.9898602708016568

$ws.Cells.Item(25, 2).Value = 1.266824619137935
$ws.Cells.Item(25, 3).Value = 0.376981551582503
$ws.Cells.Item(25, 4).Value = 0.2828236429634501
$ws.Cells.Item(25, 6).Value = 1.512447963139905
$ws.Cells.Item(25, 7).Value = 0.838982128427844
$ws.Cells.Item(25, 8).Value = 0.8470342326250488
$ws.Cells.Item(25, 9).Value = 0.5721482055093503
$ws.Cells.Item(25, 10).Value = 0.3660239290496605
$ws.Cells.Item(25, 14).Value = 1.010472419256288
